$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Apply the TableGrid table style (adds <w:tblStyle w:val="TableGrid"/> to tblPr)
$t.Style = "TableGrid"

# 2. Update the numeric values in the "t" (col 12) and "t_max" (col 13) columns
#    for each data row (rows 2-7; row 1 is the header row).
$t.Cell(2, 12).Range.Text = "0.0030"
$t.Cell(2, 13).Range.Text = "0.0131"

$t.Cell(3, 12).Range.Text = "0.0034"
$t.Cell(3, 13).Range.Text = "0.0345"

$t.Cell(4, 12).Range.Text = "0.0030"
$t.Cell(4, 13).Range.Text = "0.0282"

$t.Cell(5, 12).Range.Text = "0.0029"
$t.Cell(5, 13).Range.Text = "0.0094"

$t.Cell(6, 12).Range.Text = "0.0029"
$t.Cell(6, 13).Range.Text = "0.0082"

$t.Cell(7, 12).Range.Text = "0.0029"
$t.Cell(7, 13).Range.Text = "0.0130"

# 3. Append a new, entirely empty row (13 empty cells) at the bottom of the table.
$newRow = $t.Rows.Add()
